$d = $word.ActiveDocument
$d.Content.Find.Execute("maduras paara a ES;", $false, $false, $false, $false, $false, $true, 1, $false, "maduras para a ES;", 2)
